$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the two "footprint" bullet paragraphs (originally):
#   P1 (home)  : "“Reducing footprint at ” + “home”."
#   P2 (travel): "“reducing footprint in ” + “travel”."
# ---------------------------------------------------------------------------

$homeParaIdx = 0
$travelParaIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -clike '*Reducing footprint at home*') {
        $homeParaIdx = $i
    }
    if ($t -clike '*reducing footprint in travel*') {
        $travelParaIdx = $i
    }
}

$homePara = $d.Paragraphs.Item($homeParaIdx)

# ---------------------------------------------------------------------------
# Step 1: split the travel paragraph's leading run ("“reducing footprint in ")
# into three distinct runs: "“", "Reducing", " footprint in " -- leaving the
# trailing "travel”." run untouched. Do this first (paragraph count doesn't
# change), using InsertXML so the runtime doesn't coalesce same-format runs.
# ---------------------------------------------------------------------------

$travelPara = $d.Paragraphs.Item($travelParaIdx)
$tStart = $travelPara.Range.Start

# the leading run spans the 23 characters: “reducing footprint in  (incl. trailing space)
$oldLead = $d.Range($tStart, $tStart + 23)
if ($oldLead.Text -cne [char]8220 + "reducing footprint in ") {
    throw "unexpected travel paragraph lead text: [$($oldLead.Text)]"
}
$oldLead.Text = ""

$pkgTemplate = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>{0}</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$quote = [char]8220
$run1 = '<w:r><w:t>' + $quote + '</w:t></w:r>'
$run2 = '<w:r><w:t>Reducing</w:t></w:r>'
$run3 = '<w:r><w:t xml:space="preserve"> footprint in </w:t></w:r>'
$xmlFrag = [string]::Format($pkgTemplate, $run1 + $run2 + $run3)

$insPoint = $d.Range($tStart, $tStart)
$insPoint.InsertXML($xmlFrag)

# ---------------------------------------------------------------------------
# Step 2: insert a brand-new paragraph directly after the "home" paragraph
# carrying the old merged text, then shrink the original "home" paragraph
# down to just the quoted word "Home".
# ---------------------------------------------------------------------------

$homePara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($homeParaIdx + 1)
$newStart = $newPara.Range.Start
$newEnd = $newPara.Range.End
$newBody = $d.Range($newStart, $newEnd - 1)
$newBody.Text = $quote + "Reducing footprint at home" + [char]8221 + "."

$homeBody = $d.Range($homePara.Range.Start, $homePara.Range.End - 1)
$homeBody.Text = $quote + "Home" + [char]8221
